# Apply corrected Diebold-Mariano values to the P_valores and Estadisticos_DM sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "P_valores" ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.7422951078380673
$wsP.Range("D2").Value = 0.9136280624895987
$wsP.Range("E2").Value = 0.6728620539146042
$wsP.Range("F2").Value = 0.9444989496619038

$wsP.Range("B3").Value = 0.7422951078380673
$wsP.Range("D3").Value = 0.7722778098247107
$wsP.Range("E3").Value = 0.9238677021352635
$wsP.Range("F3").Value = 0.6721377686743928

$wsP.Range("B4").Value = 0.9136280624895987
$wsP.Range("C4").Value = 0.7722778098247107
$wsP.Range("E4").Value = 0.717555800023336
$wsP.Range("F4").Value = 0.9340578331011695

$wsP.Range("B5").Value = 0.6728620539146042
$wsP.Range("C5").Value = 0.9238677021352635
$wsP.Range("D5").Value = 0.717555800023336
$wsP.Range("F5").Value = 0.6976947698769416

$wsP.Range("B6").Value = 0.9444989496619038
$wsP.Range("C6").Value = 0.6721377686743928
$wsP.Range("D6").Value = 0.9340578331011695
$wsP.Range("E6").Value = 0.6976947698769416

# --- Sheet "Estadisticos_DM" ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = 0.3354016704513273
$wsE.Range("D2").Value = 0.1104399346997865
$wsE.Range("E2").Value = 0.4312373321681822
$wsE.Range("F2").Value = 0.07087600677957812

$wsE.Range("B3").Value = -0.3354016704513273
$wsE.Range("D3").Value = -0.2950608527102469
$wsE.Range("E3").Value = -0.09729961963251321
$wsE.Range("F3").Value = -0.4322578037199336

$wsE.Range("B4").Value = -0.1104399346997865
$wsE.Range("C4").Value = 0.2950608527102469
$wsE.Range("E4").Value = 0.3691271500125996
$wsE.Range("F4").Value = -0.08424067640261293

$wsE.Range("B5").Value = -0.4312373321681822
$wsE.Range("C5").Value = 0.09729961963251321
$wsE.Range("D5").Value = -0.3691271500125996
$wsE.Range("F5").Value = -0.3965266435774599

$wsE.Range("B6").Value = -0.07087600677957812
$wsE.Range("C6").Value = 0.4322578037199336
$wsE.Range("D6").Value = 0.08424067640261293
$wsE.Range("E6").Value = 0.3965266435774599
